$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which stores
# prices as literal strings (e.g. "218.70", "0.0622") rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.805.01"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.638.92"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "218.70"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.867.41"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "1.641.40"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "0.525"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "64.81"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "26.794.21"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "214.88"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "1.01"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "6.54"
$ws.Range("E22").Value = "  +4.09%  "
$ws.Range("D23").Value = "2.35"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "147.51"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "0.119"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "7.04"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "15.70"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "1.54"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").Value = "1.262.22"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "0.0175"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "0.528"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "0.814"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "5.33"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "1.778.09"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").Value = "92.13"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "60.17"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").Value = "1.57"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").Value = "0.0515"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "0.0962"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").Value = "7.53"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  -0.17%  "
